$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jaime Ribeiro Filho", 575, 320, 200, 1060, 955, 865, 1155, 200, 5330, 14.89),
    @("Márcio Flávio Moura de Araújo", 285, 570, 315, 385, 475, 260, 380, 80, 2750, 14.29),
    @("Fabio Miyajima", 350, 270, 330, 270, 530, 770, 80, 0, 2600, 5.26),
    @("Anya Pimentel Gomes Fernandes Vieira Meyer", 100, 310, 230, 355, 755, 300, 445, 100, 2595, 14.58),
    @("Roberto Wagner Júnior Freire de Freitas", 175, 390, 320, 365, 580, 80, 190, 20, 2120, 10),
    @("Roberto Nicolete", 160, 80, 55, 270, 335, 510, 285, 80, 1775, 15.15),
    @("Ivana Cristina de Holanda Cunha Barreto", 335, 180, 65, 150, 370, 340, 200, 0, 1640, 5.88),
    @("Caroline Pereira Bittencourt Passaes", 90, 90, 110, 360, 490, 170, 90, 90, 1490, 0),
    @("Carla Freire Celedonio Fernandes", 410, 170, 100, 220, 320, 60, 15, 170, 1465, 48),
    @("João Hermínio Martins da Silva", 310, 0, 210, 120, 55, 520, 200, 0, 1415, 0),
    @("Luiz Odorico Monteiro de Andrade", 170, 130, 15, 150, 420, 260, 140, 40, 1325, 3.7),
    @("José Luís Passos Cordeiro", 0, 400, 100, 375, 90, 130, 80, 15, 1190, 8.7),
    @("Sharmênia de Araújo Soares Nuto", 80, 55, 50, 240, 385, 115, 170, 60, 1155, 0),
    @("Fernando Ferreira Carneiro", 170, 205, 40, 130, 130, 100, 180, 180, 1135, 8.7),
    @("Marcos Roberto Lourenzoni", 110, 90, 180, 190, 120, 210, 100, 0, 1000, 20),
    @("Fernando Braga Stehling Dias", 90, 90, 90, 80, 260, 270, 80, 0, 960, 9.09),
    @("Gilvan Pessoa Furtado", 100, 180, 110, 0, 270, 100, 160, 0, 920, 14.29),
    @("Clarissa Romero Teixeira", 0, 320, 0, 80, 0, 80, 180, 80, 740, 22.22),
    @("Maximiliano Ponte", 270, 20, 75, 130, 90, 0, 95, 0, 680, 5.26),
    @("Vanira Matos Pessoa", 150, 80, 40, 20, 125, 110, 150, 0, 675, 11.11),
    @("Marcela Helena Gambim Fonseca", 0, 90, 20, 0, 170, 170, 180, 10, 640, 14.29),
    @("Regis Bernardo Brandim Gomes", 0, 0, 0, 40, 20, 160, 180, 200, 600, 9.09),
    @("Eduardo Ruback dos Santos", 0, 0, 0, 0, 0, 270, 0, 90, 360, 25),
    @("Adriana Costa Bacelo", 50, 0, 110, 0, 10, 35, 0, 80, 285, 30),
    @("Anna Carolina Machado Marinho", 0, 15, 0, 0, 150, 0, 15, 90, 270, 20),
    @("Venúcia Bruna Magalhães Pereira", 80, 80, 60, 0, 40, 0, 0, 0, 260, 0),
    @("Claudia Stutz Zubieta", 0, 0, 0, 0, 0, 90, 90, 80, 260, 0),
    @("Marlos de Medeiros Chaves", 0, 0, 0, 40, 80, 90, 0, 0, 210, 20),
    @("Galba Freire Moita", 0, 20, 95, 40, 25, 5, 0, 0, 185, 10),
    @("Alice Paula Di Sabatino Guimarães", 0, 0, 0, 0, 170, 0, 0, 0, 170, 50),
    @("Dayane Alves Costa", 90, 0, 0, 0, 40, 0, 20, 0, 150, 100),
    @("Ana Camila Oliveira Alves", 0, 0, 0, 0, 60, 60, 0, 0, 120, 50),
    @("Ana Cláudia de Araújo Teixeira", 0, 0, 0, 0, 0, 20, 90, 0, 110, 50),
    @("Margareth Borges Coutinho Gallo", 0, 0, 0, 0, 100, 0, 0, 0, 100, 25),
    @("Giovanny Augusto Camacho Antevere Mazzarotto", 80, 0, 0, 0, 5, 0, 15, 0, 100, 0),
    @("Antonio Marcos Aires Barbosa", 0, 15, 0, 0, 0, 0, 80, 0, 95, 100),
    @("Raphael Trevizani", 90, 0, 0, 0, 0, 0, 0, 0, 90, 0),
    @("Donat Alexander de Chapeaurouge", 0, 80, 0, 0, 0, 0, 0, 0, 80, 33.33)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

# Row 39 is brand new (the sheet previously ended at row 38); give column A
# the same label formatting (bold font + border + center/top alignment) used
# by every other row in column A, by copying the format from the row above.
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)

